# Auto-generated edit script: update cached Universalis market-price
# columns (currentAveragePrice*, LevePrice*, LeveProfit*) for the rows
# that changed across the ALC/ARM/BSM/CRP/CUL/GSM/WVR sheets.
$wb = $excel.ActiveWorkbook

# ALC row 26
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(26, 8).Value = 27000
$ws.Cells.Item(26, 10).Value = 27000
$ws.Cells.Item(26, 12).Value = 27000
$ws.Cells.Item(26, 14).Value = -27688

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 205.26315
$ws.Cells.Item(33, 9).Value = 142.85715
$ws.Cells.Item(33, 10).Value = 380
$ws.Cells.Item(33, 11).Value = 142.85715
$ws.Cells.Item(33, 12).Value = 380
$ws.Cells.Item(33, 13).Value = 86.14285000000001
$ws.Cells.Item(33, 14).Value = -838

# ALC row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 915.9
$ws.Cells.Item(70, 9).Value = 965.9167
$ws.Cells.Item(70, 10).Value = 903.3958
$ws.Cells.Item(70, 11).Value = 2897.7501
$ws.Cells.Item(70, 12).Value = 2710.1874
$ws.Cells.Item(70, 13).Value = -2627.7501
$ws.Cells.Item(70, 14).Value = -3250.1874

# ALC row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73, 8).Value = 915.9
$ws.Cells.Item(73, 9).Value = 965.9167
$ws.Cells.Item(73, 10).Value = 903.3958
$ws.Cells.Item(73, 11).Value = 2897.7501
$ws.Cells.Item(73, 12).Value = 2710.1874
$ws.Cells.Item(73, 13).Value = -1961.7501
$ws.Cells.Item(73, 14).Value = -4582.1874

# ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(88, 8).Value = 1398.1666
$ws.Cells.Item(88, 9).Value = 1126
$ws.Cells.Item(88, 10).Value = 1502.8462
$ws.Cells.Item(88, 11).Value = 1126
$ws.Cells.Item(88, 12).Value = 1502.8462
$ws.Cells.Item(88, 13).Value = -720
$ws.Cells.Item(88, 14).Value = -2314.8462

# ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(91, 8).Value = 1398.1666
$ws.Cells.Item(91, 9).Value = 1126
$ws.Cells.Item(91, 10).Value = 1502.8462
$ws.Cells.Item(91, 11).Value = 1126
$ws.Cells.Item(91, 12).Value = 1502.8462
$ws.Cells.Item(91, 13).Value = 278
$ws.Cells.Item(91, 14).Value = -4310.8462

# ALC row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(92, 8).Value = 6944837.5
$ws.Cells.Item(92, 9).Value = 8064789
$ws.Cells.Item(92, 10).Value = 1140.2
$ws.Cells.Item(92, 11).Value = 8064789
$ws.Cells.Item(92, 12).Value = 1140.2
$ws.Cells.Item(92, 13).Value = -8063541
$ws.Cells.Item(92, 14).Value = -3636.2

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1536.1489
$ws.Cells.Item(137, 9).Value = 1180.6129
$ws.Cells.Item(137, 10).Value = 2225
$ws.Cells.Item(137, 11).Value = 3541.8387
$ws.Cells.Item(137, 12).Value = 6675
$ws.Cells.Item(137, 13).Value = -991.8387000000002
$ws.Cells.Item(137, 14).Value = -11775

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2077.5146
$ws.Cells.Item(138, 9).Value = 1609.0454
$ws.Cells.Item(138, 10).Value = 2301.5652
$ws.Cells.Item(138, 11).Value = 4827.1362
$ws.Cells.Item(138, 12).Value = 6904.6956
$ws.Cells.Item(138, 13).Value = 312.8638000000001
$ws.Cells.Item(138, 14).Value = -17184.6956

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2429.2632
$ws.Cells.Item(2, 9).Value = 1610.0769
$ws.Cells.Item(2, 10).Value = 4204.1665
$ws.Cells.Item(2, 11).Value = 1610.0769
$ws.Cells.Item(2, 12).Value = 4204.1665
$ws.Cells.Item(2, 13).Value = -1497.0769
$ws.Cells.Item(2, 14).Value = -4430.1665

# ARM row 36
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(36, 8).Value = 6099.2
$ws.Cells.Item(36, 9).Value = 1498
$ws.Cells.Item(36, 10).Value = 9166.666999999999
$ws.Cells.Item(36, 11).Value = 1498
$ws.Cells.Item(36, 12).Value = 9166.666999999999
$ws.Cells.Item(36, 13).Value = -1152
$ws.Cells.Item(36, 14).Value = -9858.666999999999

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 2953.158
$ws.Cells.Item(63, 9).Value = 1821
$ws.Cells.Item(63, 10).Value = 4211.1113
$ws.Cells.Item(63, 11).Value = 1821
$ws.Cells.Item(63, 12).Value = 4211.1113
$ws.Cells.Item(63, 13).Value = -1135
$ws.Cells.Item(63, 14).Value = -5583.1113

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(66, 8).Value = 2953.158
$ws.Cells.Item(66, 9).Value = 1821
$ws.Cells.Item(66, 10).Value = 4211.1113
$ws.Cells.Item(66, 11).Value = 9105
$ws.Cells.Item(66, 12).Value = 21055.5565
$ws.Cells.Item(66, 13).Value = -5673
$ws.Cells.Item(66, 14).Value = -27919.5565

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 2429.2632
$ws.Cells.Item(116, 9).Value = 1610.0769
$ws.Cells.Item(116, 10).Value = 4204.1665
$ws.Cells.Item(116, 11).Value = 1610.0769
$ws.Cells.Item(116, 12).Value = 4204.1665
$ws.Cells.Item(116, 13).Value = 683.9231
$ws.Cells.Item(116, 14).Value = -8792.166499999999

# ARM row 123
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(123, 8).Value = 355000
$ws.Cells.Item(123, 10).Value = 355000
$ws.Cells.Item(123, 12).Value = 355000
$ws.Cells.Item(123, 14).Value = -364800

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2429.2632
$ws.Cells.Item(3, 9).Value = 1610.0769
$ws.Cells.Item(3, 10).Value = 4204.1665
$ws.Cells.Item(3, 11).Value = 1610.0769
$ws.Cells.Item(3, 12).Value = 4204.1665
$ws.Cells.Item(3, 13).Value = -1496.0769
$ws.Cells.Item(3, 14).Value = -4432.1665

# BSM row 64
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 1043
$ws.Cells.Item(64, 9).Value = 975.5
$ws.Cells.Item(64, 10).Value = 1133
$ws.Cells.Item(64, 11).Value = 975.5
$ws.Cells.Item(64, 12).Value = 1133
$ws.Cells.Item(64, 13).Value = -750.5
$ws.Cells.Item(64, 14).Value = -1583

# BSM row 67
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(67, 8).Value = 1043
$ws.Cells.Item(67, 9).Value = 975.5
$ws.Cells.Item(67, 10).Value = 1133
$ws.Cells.Item(67, 11).Value = 975.5
$ws.Cells.Item(67, 12).Value = 1133
$ws.Cells.Item(67, 13).Value = -195.5
$ws.Cells.Item(67, 14).Value = -2693

# CRP row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 33.357143
$ws.Cells.Item(7, 9).Value = 15
$ws.Cells.Item(7, 10).Value = 66.40000000000001
$ws.Cells.Item(7, 11).Value = 15
$ws.Cells.Item(7, 12).Value = 66.40000000000001
$ws.Cells.Item(7, 13).Value = 98
$ws.Cells.Item(7, 14).Value = -292.4

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 3448.6667
$ws.Cells.Item(16, 9).Value = 3448.6667
$ws.Cells.Item(16, 11).Value = 3448.6667
$ws.Cells.Item(16, 13).Value = -3161.6667

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2314.6897
$ws.Cells.Item(58, 9).Value = 2437.5454
$ws.Cells.Item(58, 10).Value = 1928.5714
$ws.Cells.Item(58, 11).Value = 2437.5454
$ws.Cells.Item(58, 12).Value = 1928.5714
$ws.Cells.Item(58, 13).Value = -2234.5454
$ws.Cells.Item(58, 14).Value = -2334.5714

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 5050.846
$ws.Cells.Item(105, 9).Value = 4697.625
$ws.Cells.Item(105, 10).Value = 5616
$ws.Cells.Item(105, 11).Value = 4697.625
$ws.Cells.Item(105, 12).Value = 5616
$ws.Cells.Item(105, 13).Value = -2950.625
$ws.Cells.Item(105, 14).Value = -9110

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113, 8).Value = 3448.6667
$ws.Cells.Item(113, 9).Value = 3448.6667
$ws.Cells.Item(113, 11).Value = 3448.6667
$ws.Cells.Item(113, 13).Value = -1278.6667

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 2314.6897
$ws.Cells.Item(136, 9).Value = 2437.5454
$ws.Cells.Item(136, 10).Value = 1928.5714
$ws.Cells.Item(136, 11).Value = 7312.6362
$ws.Cells.Item(136, 12).Value = 5785.7142
$ws.Cells.Item(136, 13).Value = -4762.6362
$ws.Cells.Item(136, 14).Value = -10885.7142

# CUL row 19
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(19, 8).Value = 1100
$ws.Cells.Item(19, 10).Value = 1100
$ws.Cells.Item(19, 12).Value = 3300
$ws.Cells.Item(19, 14).Value = -3648

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 3637.1785
$ws.Cells.Item(122, 10).Value = 3840.1646
$ws.Cells.Item(122, 12).Value = 34561.4814
$ws.Cells.Item(122, 14).Value = -39461.4814

# GSM row 109
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(109, 8).Value = 0
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 14).ClearContents()

# WVR row 31
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 14).ClearContents()
